# Export latest mex files to xlsx. Update segments and articles_db.
# Append 10 new coded-segment rows (379-388) to Sheet1, mirroring the
# existing table's layout/formatting (row 378 is used as the style
# template for every new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# ---------------------------------------------------------------------
# 1) Stamp every new row (379-388) with the same cell formatting
#    (fills/borders/fonts/alignment) as the last existing data row, so
#    styles line up (s="6","2","1","2","1","1","1","3","2","3","4","1","1").
# ---------------------------------------------------------------------
$ws.Range("A378:M378").Copy()
$ws.Range("A379:M388").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row data (mirrors the rows added in the source diff).
# ---------------------------------------------------------------------
$rows = @(
  @{ Row=379; E="Bacteria:Binomial (genus species)"; F="1: 709";  G="1: 722";  H=0; I=".  `npneumoniae"; J=13; K=0.040218; M="11/8/18 14:41:00" },
  @{ Row=380; E="Bacteria:Binomial (genus species)"; F="1: 727";  G="1: 740";  H=0; I=" H. influenzae";  J=13; K=0.040218; M="11/8/18 14:41:00" },
  @{ Row=381; E="Event month";                       F="2: 2174"; G="2: 2180"; H=0; I="October";         J=7;  K=0.021656; M="11/13/18 08:38:00" },
  @{ Row=382; E="Event month";                       F="2: 2192"; G="2: 2196"; H=0; I="March";           J=5;  K=0.015468; M="11/13/18 08:38:00" },
  @{ Row=383; E="Event year";                        F="2: 2182"; G="2: 2185"; H=0; I="2004";            J=4;  K=0.012375; M="11/13/18 08:38:00" },
  @{ Row=384; E="Event year";                        F="2: 2198"; G="2: 2201"; H=0; I="2006";            J=4;  K=0.012375; M="11/13/18 08:38:00" },
  @{ Row=385; E="B";                                  F="2: 2198"; G="2: 2201"; H=0; I="2006";            J=4;  K=0.012375; M="11/13/18 08:38:00" },
  @{ Row=386; E="B";                                  F="2: 2192"; G="2: 2196"; H=0; I="March";           J=5;  K=0.015468; M="11/13/18 08:38:00" },
  @{ Row=387; E="B";                                  F="2: 2174"; G="2: 2180"; H=0; I="October";         J=7;  K=0.021656; M="11/13/18 08:38:00" },
  @{ Row=388; E="B";                                  F="2: 2182"; G="2: 2185"; H=0; I="2004";            J=4;  K=0.012375; M="11/13/18 08:38:00" }
)

$blackCircle = [char]0x25CF

foreach ($r in $rows) {
    $row = $r.Row

    # Color / Document group / Document name are identical for every new
    # row: the "bullet" marker, blank comment/doc-group, and the
    # "18833" document id (already an existing shared string).
    $ws.Cells.Item($row, 1).Value = $blackCircle
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = "Sonia"
    $ws.Cells.Item($row, 13).Value = $r.M

    # Segment text (column I) - write directly unless it is purely
    # numeric-looking ("2004"/"2006"), in which case force text below.
    $looksNumeric = $r.I -match '^[0-9]+$'

    if ($looksNumeric) {
        $ws.Cells.Item($row, 9).NumberFormat = "@"
        $ws.Cells.Item($row, 9).Value = $r.I
        # Re-apply the template's number format/style so the cell keeps
        # style index 2 (General) instead of the ad-hoc "@" format, while
        # preserving the text we just stored.
        $ws.Range("I378").Copy()
        $ws.Range("I" + $row).PasteSpecial($xlPasteFormats)
        $excel.CutCopyMode = 0
    } else {
        $ws.Cells.Item($row, 9).Value = $r.I
    }
}

# ---------------------------------------------------------------------
# 3) Column D ("18833") is a numeric-looking Document name that must stay
#    text (it already exists as a shared string elsewhere, e.g. D5), so
#    copy its *value* (not its format) into every new row's D cell.
# ---------------------------------------------------------------------
$ws.Range("D5").Copy()
$ws.Range("D379:D388").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Columns B and C stay blank (same empty string used throughout the
#    sheet) - the format copy in step 1 already created those cells, so
#    nothing further is required there.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5) Row heights: row 379 wraps the two-line "pneumoniae" segment (30pt,
#    i.e. 2 lines), the rest are single-line rows (16pt), matching the
#    rest of the table.
# ---------------------------------------------------------------------
$ws.Rows.Item(379).RowHeight = 30
$ws.Range("380:388").RowHeight = 16
